$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("F2").Value = "AnswZ"
$ws.Range("G2").Value = "PruebaZ"
$ws.Range("H2").Value = 20300120

$ws.Range("F3").Value = "MattioliX"
$ws.Range("G3").Value = "PruebaX"
$ws.Range("H3").Value = 20300121

$ws.Range("F5").Select()
